$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.3419937204270411 ; $ws.Cells.Item(2, 3).Value = 0.07295521617012923 ; $ws.Cells.Item(2, 4).Value = 0.0790759291399894 ; $ws.Cells.Item(2, 5).Value = 0.4186173600666478 ; $ws.Cells.Item(2, 7).Value = 0.4099646853786965 ; $ws.Cells.Item(2, 8).Value = 0.566722113727046 ; $ws.Cells.Item(2, 9).Value = 0.431979201775178 ; $ws.Cells.Item(2, 11).Value = 0.3789283569079771 ; $ws.Cells.Item(2, 15).Value = 1.900787633321301
$ws.Cells.Item(3, 2).Value = 0.2994647181646997 ; $ws.Cells.Item(3, 3).Value = 0.06502746772761725 ; $ws.Cells.Item(3, 4).Value = 0.07165833778176989 ; $ws.Cells.Item(3, 5).Value = 0.365274963597173 ; $ws.Cells.Item(3, 7).Value = 0.4105030486526928 ; $ws.Cells.Item(3, 8).Value = 0.5710512420052467 ; $ws.Cells.Item(3, 9).Value = 0.437165556847571 ; $ws.Cells.Item(3, 11).Value = 0.3305766105088139 ; $ws.Cells.Item(3, 15).Value = 1.910712916604183
$ws.Cells.Item(4, 2).Value = 0.2732938709733048 ; $ws.Cells.Item(4, 3).Value = 0.06013130166671488 ; $ws.Cells.Item(4, 4).Value = 0.06713703597539222 ; $ws.Cells.Item(4, 5).Value = 0.3326033741979728 ; $ws.Cells.Item(4, 7).Value = 0.4111944143402155 ; $ws.Cells.Item(4, 8).Value = 0.5740103714446647 ; $ws.Cells.Item(4, 9).Value = 0.4406458533377595 ; $ws.Cells.Item(4, 11).Value = 0.3007926697861762 ; $ws.Cells.Item(4, 15).Value = 1.918192594569561
$ws.Cells.Item(5, 2).Value = 0.2626150629347421 ; $ws.Cells.Item(5, 3).Value = 0.05812897051117716 ; $ws.Cells.Item(5, 4).Value = 0.06530291569288238 ; $ws.Cells.Item(5, 5).Value = 0.3193082460111043 ; $ws.Cells.Item(5, 7).Value = 0.4115666001412137 ; $ws.Cells.Item(5, 8).Value = 0.5752918947514729 ; $ws.Cells.Item(5, 9).Value = 0.4421383703937352 ; $ws.Cells.Item(5, 11).Value = 0.2886319761973368 ; $ws.Cells.Item(5, 15).Value = 1.921588336167574
$ws.Cells.Item(6, 2).Value = 0.2608410279604243 ; $ws.Cells.Item(6, 3).Value = 0.05779605774353058 ; $ws.Cells.Item(6, 4).Value = 0.06499886603796767 ; $ws.Cells.Item(6, 5).Value = 0.3171016918152105 ; $ws.Cells.Item(6, 7).Value = 0.411633856004471 ; $ws.Cells.Item(6, 8).Value = 0.5755092588289443 ; $ws.Cells.Item(6, 9).Value = 0.442390683606444 ; $ws.Cells.Item(6, 11).Value = 0.2866113005742079 ; $ws.Cells.Item(6, 15).Value = 1.922173177353173
$ws.Cells.Item(7, 2).Value = 0.2731499086519875 ; $ws.Cells.Item(7, 3).Value = 0.06010432620882966 ; $ws.Cells.Item(7, 4).Value = 0.06711226656943836 ; $ws.Cells.Item(7, 5).Value = 0.3324239972883305 ; $ws.Cells.Item(7, 7).Value = 0.4111990679172308 ; $ws.Cells.Item(7, 8).Value = 0.5740273482536864 ; $ws.Cells.Item(7, 9).Value = 0.4406656814063155 ; $ws.Cells.Item(7, 11).Value = 0.3006287607718434 ; $ws.Cells.Item(7, 15).Value = 1.918236983773483
$ws.Cells.Item(8, 2).Value = 0.3273420625438348 ; $ws.Cells.Item(8, 3).Value = 0.07022767759735871 ; $ws.Cells.Item(8, 4).Value = 0.07651147508599365 ; $ws.Cells.Item(8, 5).Value = 0.4002071868227546 ; $ws.Cells.Item(8, 7).Value = 0.4100752578014024 ; $ws.Cells.Item(8, 8).Value = 0.568152293517393 ; $ws.Cells.Item(8, 9).Value = 0.4337059764904616 ; $ws.Cells.Item(8, 11).Value = 0.3622768889098609 ; $ws.Cells.Item(8, 15).Value = 1.903921891058133
$ws.Cells.Item(9, 2).Value = 0.433133462145264 ; $ws.Cells.Item(9, 3).Value = 0.08985165565974285 ; $ws.Cells.Item(9, 4).Value = 0.09520643689737085 ; $ws.Cells.Item(9, 5).Value = 0.5338462420115491 ; $ws.Cells.Item(9, 7).Value = 0.4107475344981282 ; $ws.Cells.Item(9, 8).Value = 0.5590219255403497 ; $ws.Cells.Item(9, 9).Value = 0.4224111945678963 ; $ws.Cells.Item(9, 11).Value = 0.4823888465412836 ; $ws.Cells.Item(9, 15).Value = 1.886876139744686
$ws.Cells.Item(10, 2).Value = 0.5105469197237085 ; $ws.Cells.Item(10, 3).Value = 0.1041294763322753 ; $ws.Cells.Item(10, 4).Value = 0.1091038841048402 ; $ws.Cells.Item(10, 5).Value = 0.6325801516627507 ; $ws.Cells.Item(10, 7).Value = 0.4130139381484526 ; $ws.Cells.Item(10, 8).Value = 0.5537743890861861 ; $ws.Cells.Item(10, 9).Value = 0.415555380923557 ; $ws.Cells.Item(10, 11).Value = 0.570141579896756 ; $ws.Cells.Item(10, 15).Value = 1.881122114504421
$ws.Cells.Item(11, 2).Value = 0.5456930962515969 ; $ws.Cells.Item(11, 3).Value = 0.1105943320977758 ; $ws.Cells.Item(11, 4).Value = 0.1154619028305461 ; $ws.Cells.Item(11, 5).Value = 0.6776409264464149 ; $ws.Cells.Item(11, 7).Value = 0.4144340397891995 ; $ws.Cells.Item(11, 8).Value = 0.5517050126349829 ; $ws.Cells.Item(11, 9).Value = 0.4127514227220495 ; $ws.Cells.Item(11, 11).Value = 0.6099521553671252 ; $ws.Cells.Item(11, 15).Value = 1.879984846479402
$ws.Cells.Item(12, 2).Value = 0.5589915640219658 ; $ws.Cells.Item(12, 3).Value = 0.1130380222911924 ; $ws.Cells.Item(12, 4).Value = 0.1178747007182892 ; $ws.Cells.Item(12, 5).Value = 0.6947271311692589 ; $ws.Cells.Item(12, 7).Value = 0.4150280705964633 ; $ws.Cells.Item(12, 8).Value = 0.5509671449164983 ; $ws.Cells.Item(12, 9).Value = 0.4117350464549325 ; $ws.Cells.Item(12, 11).Value = 0.625011307833546 ; $ws.Cells.Item(12, 15).Value = 1.879767866987891
$ws.Cells.Item(13, 2).Value = 0.5561279834765855 ; $ws.Cells.Item(13, 3).Value = 0.1125119271140989 ; $ws.Cells.Item(13, 4).Value = 0.1173548326529783 ; $ws.Cells.Item(13, 5).Value = 0.6910462768777563 ; $ws.Cells.Item(13, 7).Value = 0.4148976278063685 ; $ws.Cells.Item(13, 8).Value = 0.5511240216500681 ; $ws.Cells.Item(13, 9).Value = 0.411951918632429 ; $ws.Cells.Item(13, 11).Value = 0.6217687826131453 ; $ws.Cells.Item(13, 15).Value = 1.879805080909335
$ws.Cells.Item(14, 2).Value = 0.5467873857009522 ; $ws.Cells.Item(14, 3).Value = 0.1107954647394251 ; $ws.Cells.Item(14, 4).Value = 0.1156603019985027 ; $ws.Cells.Item(14, 5).Value = 0.679046155061684 ; $ws.Cells.Item(14, 7).Value = 0.4144817815305828 ; $ws.Cells.Item(14, 8).Value = 0.5516433903215301 ; $ws.Cells.Item(14, 9).Value = 0.4126668936187095 ; $ws.Cells.Item(14, 11).Value = 0.6111914101640252 ; $ws.Cells.Item(14, 15).Value = 1.879962708302145
$ws.Cells.Item(15, 2).Value = 0.5410645976724311 ; $ws.Cells.Item(15, 3).Value = 0.1097435056013865 ; $ws.Cells.Item(15, 4).Value = 0.1146230235927277 ; $ws.Cells.Item(15, 5).Value = 0.6716987348148677 ; $ws.Cells.Item(15, 7).Value = 0.4142344009781453 ; $ws.Cells.Item(15, 8).Value = 0.5519674799778898 ; $ws.Cells.Item(15, 9).Value = 0.4131107563433538 ; $ws.Cells.Item(15, 11).Value = 0.6047103293468012 ; $ws.Cells.Item(15, 15).Value = 1.880087111294557
$ws.Cells.Item(16, 2).Value = 0.5082485810177673 ; $ws.Cells.Item(16, 3).Value = 0.10370636932754 ; $ws.Cells.Item(16, 4).Value = 0.1086890951216049 ; $ws.Cells.Item(16, 5).Value = 0.6296384063867464 ; $ws.Cells.Item(16, 7).Value = 0.4129289879773808 ; $ws.Cells.Item(16, 8).Value = 0.5539160274922779 ; $ws.Cells.Item(16, 9).Value = 0.4157449738133714 ; $ws.Cells.Item(16, 11).Value = 0.5675376258630536 ; $ws.Cells.Item(16, 15).Value = 1.881226294346334
$ws.Cells.Item(17, 2).Value = 0.4880987513508046 ; $ws.Cells.Item(17, 3).Value = 0.09999499895883446 ; $ws.Cells.Item(17, 4).Value = 0.1050580267680061 ; $ws.Cells.Item(17, 5).Value = 0.6038743717326582 ; $ws.Cells.Item(17, 7).Value = 0.4122280432832639 ; $ws.Cells.Item(17, 8).Value = 0.5551928306325067 ; $ws.Cells.Item(17, 9).Value = 0.4174417160432391 ; $ws.Cells.Item(17, 11).Value = 0.5447051127510463 ; $ws.Cells.Item(17, 15).Value = 1.882304860174258
$ws.Cells.Item(18, 2).Value = 0.476502588530451 ; $ws.Cells.Item(18, 3).Value = 0.09785747867576333 ; $ws.Cells.Item(18, 4).Value = 0.1029729221839659 ; $ws.Cells.Item(18, 5).Value = 0.5890692449155779 ; $ws.Cells.Item(18, 7).Value = 0.411861483740708 ; $ws.Cells.Item(18, 8).Value = 0.5559571184720369 ; $ws.Cells.Item(18, 9).Value = 0.4184472570759823 ; $ws.Cells.Item(18, 11).Value = 0.5315622723816205 ; $ws.Cells.Item(18, 15).Value = 1.883064503384333
$ws.Cells.Item(19, 2).Value = 0.4725752272387354 ; $ws.Cells.Item(19, 3).Value = 0.09713326566013336 ; $ws.Cells.Item(19, 4).Value = 0.1022675244988847 ; $ws.Cells.Item(19, 5).Value = 0.5840587804662363 ; $ws.Cells.Item(19, 7).Value = 0.4117436503326672 ; $ws.Cells.Item(19, 8).Value = 0.5562210272735513 ; $ws.Cells.Item(19, 9).Value = 0.418792798470502 ; $ws.Cells.Item(19, 11).Value = 0.52711060296852 ; $ws.Cells.Item(19, 15).Value = 1.883345602203718
$ws.Cells.Item(20, 2).Value = 0.4902444137623263 ; $ws.Cells.Item(20, 3).Value = 0.1003903749126778 ; $ws.Cells.Item(20, 4).Value = 0.1054442095149426 ; $ws.Cells.Item(20, 5).Value = 0.6066155711565813 ; $ws.Cells.Item(20, 7).Value = 0.4122988692360394 ; $ws.Cells.Item(20, 8).Value = 0.5550538172705757 ; $ws.Cells.Item(20, 9).Value = 0.4172580282690141 ; $ws.Cells.Item(20, 11).Value = 0.547136731234616 ; $ws.Cells.Item(20, 15).Value = 1.882175624038723
$ws.Cells.Item(21, 2).Value = 0.5495312390005154 ; $ws.Cells.Item(21, 3).Value = 0.1112997515819529 ; $ws.Cells.Item(21, 4).Value = 0.1161578870272137 ; $ws.Cells.Item(21, 5).Value = 0.6825702534044922 ; $ws.Cells.Item(21, 7).Value = 0.4146023960103804 ; $ws.Cells.Item(21, 8).Value = 0.5514895966920079 ; $ws.Cells.Item(21, 9).Value = 0.4124556541530033 ; $ws.Cells.Item(21, 11).Value = 0.6142986868627531 ; $ws.Cells.Item(21, 15).Value = 1.8799106035911
$ws.Cells.Item(22, 2).Value = 0.5882163827702414 ; $ws.Cells.Item(22, 3).Value = 0.1184039217787927 ; $ws.Cells.Item(22, 4).Value = 0.1231899263840006 ; $ws.Cells.Item(22, 5).Value = 0.7323439220684094 ; $ws.Cells.Item(22, 7).Value = 0.4164360104881979 ; $ws.Cells.Item(22, 8).Value = 0.5494269299599353 ; $ws.Cells.Item(22, 9).Value = 0.4095818509717724 ; $ws.Cells.Item(22, 11).Value = 0.6580979077889992 ; $ws.Cells.Item(22, 15).Value = 1.879676148838684
$ws.Cells.Item(23, 2).Value = 0.5675752904312787 ; $ws.Cells.Item(23, 3).Value = 0.1146146697862775 ; $ws.Cells.Item(23, 4).Value = 0.1194340542940324 ; $ws.Cells.Item(23, 5).Value = 0.7057660731069575 ; $ws.Cells.Item(23, 7).Value = 0.4154272477123016 ; $ws.Cells.Item(23, 8).Value = 0.5505033832049406 ; $ws.Cells.Item(23, 9).Value = 0.4110913705984061 ; $ws.Cells.Item(23, 11).Value = 0.6347303401447562 ; $ws.Cells.Item(23, 15).Value = 1.879687013973779
$ws.Cells.Item(24, 2).Value = 0.4892743963633279 ; $ws.Cells.Item(24, 3).Value = 0.1002116372892772 ; $ws.Cells.Item(24, 4).Value = 0.1052696086791371 ; $ws.Cells.Item(24, 5).Value = 0.6053762533719862 ; $ws.Cells.Item(24, 7).Value = 0.4122667353889966 ; $ws.Cells.Item(24, 8).Value = 0.5551165710271135 ; $ws.Cells.Item(24, 9).Value = 0.4173409798814625 ; $ws.Cells.Item(24, 11).Value = 0.5460374466983353 ; $ws.Cells.Item(24, 15).Value = 1.882233617018301
$ws.Cells.Item(25, 2).Value = 0.4045674058705799 ; $ws.Cells.Item(25, 3).Value = 0.08456735182578257 ; $ws.Cells.Item(25, 4).Value = 0.09012059176609455 ; $ws.Cells.Item(25, 5).Value = 0.4976050485609704 ; $ws.Cells.Item(25, 7).Value = 0.4102557836502427 ; $ws.Cells.Item(25, 8).Value = 0.5612356766135775 ; $ws.Cells.Item(25, 9).Value = 0.4252139928083416 ; $ws.Cells.Item(25, 11).Value = 0.4499807566113816 ; $ws.Cells.Item(25, 15).Value = 1.890302124189645
